# Swap the data of row 2 and row 3 (the Tjäder / Spillkråka records),
# moving the "Publik kommentar" note ("2 tuppar") from row 2 to row 3
# along with the Tjäder record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold data for these two records.
$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}

# Move the public comment ("2 tuppar") from row 2 to row 3.
$ws.Range("AC2").Value2 = ""
$ws.Range("AC3").Value2 = "2 tuppar"
